# Updated symbol list on Thu Jan 12 11:35:01 UTC 2023 with GitHub Actions
#
# The "cryptos" sheet refreshes its scraped coinranking.com snapshot:
#   - Coin/Link (columns B/C) for rows 7-18 shift down one slot (the
#     exchange-token board churned - a new coin "GateToken" entered at the
#     top of that block and every other row's old neighbour slid down),
#     while row 7 picks up the coin that fell off the bottom (row 18's old
#     GateToken row).
#   - Price (D) and Volume(1h) (E) are refreshed for every row that still
#     had live market data.
#
# All of these columns are stored as literal text in the workbook (e.g.
# "0.06510" keeps its trailing zero, "1.72%" is the literal string, not a
# percentage-formatted number), so each cell is pre-formatted as Text
# before the new value is written - otherwise Excel would happily "helpfully"
# coerce a numeric-looking string into a real number and silently mangle
# trailing zeros / turn the "%" strings into percentage-formatted numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Coin / Link rotation across rows 7-18 ------------------------------
$coinRotation = @(
    @{ Row = 7;  Coin = "GateToken";                         Link = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Row = 8;  Coin = "FTXToken";                           Link = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Row = 9;  Coin = "MXToken";                            Link = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Row = 10; Coin = "WazirX";                             Link = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Row = 11; Coin = "LiechtensteinCryptoassetsExchange";  Link = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Row = 12; Coin = "MandalaExchangeToken";               Link = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Row = 13; Coin = "BitrueCoin";                         Link = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Row = 14; Coin = "BitMartToken";                       Link = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Row = 15; Coin = "BitForexToken";                      Link = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Row = 16; Coin = "One";                                Link = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" },
    @{ Row = 17; Coin = "TigerCash";                          Link = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Row = 18; Coin = "LEO";                                Link = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
)

foreach ($entry in $coinRotation) {
    $ws.Range("B$($entry.Row)").Value = $entry.Coin
    $ws.Range("C$($entry.Row)").Value = $entry.Link
}

# --- Refreshed Price (D) / Volume(1h) (E) values -------------------------
$priceVolume = @(
    @{ Row = 2;  Price = "282.82";     Volume = "1.72%" },
    @{ Row = 3;  Price = "28.34";      Volume = "4.02%" },
    @{ Row = 4;  Price = "5.014";      Volume = "3.14%" },
    @{ Row = 5;  Price = "0.06510";    Volume = "1.26%" },
    @{ Row = 6;  Price = "7.218";      Volume = "3.03%" },
    @{ Row = 7;  Price = "3.352";      Volume = "1.40%" },
    @{ Row = 8;  Price = "1.389";      Volume = "16.41%" },
    @{ Row = 9;  Price = "0.9182";     Volume = "3.74%" },
    @{ Row = 10; Price = "0.1536";     Volume = "-0.29%" },
    @{ Row = 11; Price = "0.06522";    Volume = "26.03%" },
    @{ Row = 12; Price = "0.07552";    Volume = "0.77%" },
    @{ Row = 13; Price = "0.02846";    Volume = "-1.19%" },
    @{ Row = 14; Price = "0.08981";    Volume = "0.01%" },
    @{ Row = 15; Price = "0.001585";   Volume = "0.97%" },
    @{ Row = 16; Price = "0.0006382";  Volume = "0.26%" },
    @{ Row = 17; Price = "0.006190";   Volume = "0.93%" },
    @{ Row = 18; Price = "3.445";      Volume = "-0.93%" },
    @{ Row = 19; Price = "2.236";      Volume = "-0.28%" },
    @{ Row = 20; Price = $null;        Volume = "-0.05%" },
    @{ Row = 21; Price = "0.1283";     Volume = "-4.38%" },
    @{ Row = 22; Price = "3.968";      Volume = "1.51%" },
    @{ Row = 23; Price = "0.1544";     Volume = "1.74%" },
    @{ Row = 24; Price = "0.04444";    Volume = "0.67%" },
    @{ Row = 25; Price = $null;        Volume = "0.73%" },
    @{ Row = 26; Price = "0.004434";   Volume = "14.08%" },
    @{ Row = 27; Price = $null;        Volume = "1.66%" },
    @{ Row = 28; Price = $null;        Volume = "-1.55%" },
    @{ Row = 40; Price = "0.04117";    Volume = "0.08%" },
    @{ Row = 41; Price = "0.006699";   Volume = "-1.88%" },
    @{ Row = 43; Price = "0.002189";   Volume = "14.62%" },
    @{ Row = 44; Price = "0.01207";    Volume = "3.51%" },
    @{ Row = 45; Price = "0.00005668"; Volume = "6.58%" },
    @{ Row = 46; Price = "1.965";      Volume = "16.72%" }
)

foreach ($entry in $priceVolume) {
    if ($null -ne $entry.Price) {
        Set-TextValue "D$($entry.Row)" $entry.Price
    }
    if ($null -ne $entry.Volume) {
        Set-TextValue "E$($entry.Row)" $entry.Volume
    }
}
